$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data row for account 005681354 / MATHEUS / 82700 (row 3:
# header is row 1, BRASFORT is row 2, MATHEUS is row 3). Deleting the
# entire row shifts all subsequent rows up by one, matching the diff.
$ws.Rows.Item(3).Delete()
